$d = $word.ActiveDocument

# --- 1) Merge runs (remove proofErr spell-check splits) by replacing each
#        paragraph's text with itself via Find/Replace, which Word's COM
#        layer re-serializes as a single run. ---

$d.Content.Find.Execute("“Next Track” onClick should update the currently playing song information", $true, $false, $false, $false, $false, $true, 1, $false, "“Next Track” onClick should update the currently playing song information", 2) | Out-Null

$d.Content.Find.Execute("In manipulateAPlaylist, addSongsFromSongBank, with an empty playlist their will be no image for it, however, when you add songs to it the image should appear. I need to make that happen.", $true, $false, $false, $false, $false, $true, 1, $false, "In manipulateAPlaylist, addSongsFromSongBank, with an empty playlist their will be no image for it, however, when you add songs to it the image should appear. I need to make that happen.", 2) | Out-Null

$d.Content.Find.Execute("Update all components to user the SpotifyAPIBaseComposition", $true, $false, $false, $false, $false, $true, 1, $false, "Update all components to user the SpotifyAPIBaseComposition", 2) | Out-Null

$d.Content.Find.Execute("Make it so on request if the access token is expired the app refreshes it and then calls the function that failed. This will increase the apps reponsivness", $true, $false, $false, $false, $false, $true, 1, $false, "Make it so on request if the access token is expired the app refreshes it and then calls the function that failed. This will increase the apps reponsivness", 2) | Out-Null

$d.Content.Find.Execute("Give everything the react api composition component", $true, $false, $false, $false, $false, $true, 1, $false, "Give everything the react api composition component", 2) | Out-Null

# --- 2) Append a tab after "Compare two playlists " (last paragraph) ---
$lastP = $d.Paragraphs($d.Paragraphs.Count)
$lastR = $lastP.Range
$lastR.Collapse(0)
$lastR.InsertAfter([char]9)

# --- 3) Add three new bullet paragraphs (ListParagraph, ilvl=0, numId=1)
#        after the "Compare two playlists" item. ---
$p1 = $d.Paragraphs($d.Paragraphs.Count)
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$newP1 = $d.Paragraphs($d.Paragraphs.Count)
$newR1 = $newP1.Range
$newR1.ListFormat.ListLevelNumber = 1
$newR1.Text = "Work on data storage that stores the users information so we do not have to keep on making requests and increases responsiveness. At least for the songs in the song bank."

$p2 = $d.Paragraphs($d.Paragraphs.Count)
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$newP2 = $d.Paragraphs($d.Paragraphs.Count)
$newR2 = $newP2.Range
$newR2.ListFormat.ListLevelNumber = 1
$newR2.Text = "Search Songs"

$p3 = $d.Paragraphs($d.Paragraphs.Count)
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$newP3 = $d.Paragraphs($d.Paragraphs.Count)
$newR3 = $newP3.Range
$newR3.ListFormat.ListLevelNumber = 1
$newR3.Text = "Delete Songs in compare two playlists"

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
